$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 16454
$ws1.Range("F5").Value = 739
$ws1.Range("F6").Value = 15619
$ws1.Range("F8").Value = 9286
$ws1.Range("F9").Value = 500
$ws1.Range("G9").Value = 70
$ws1.Range("G10").Value = 55
$ws1.Range("F11").Value = 1035
$ws1.Range("F12").Value = 130
$ws1.Range("G13").Value = 50
$ws1.Range("G15").Value = 60
$ws1.Range("F16").Value = 23
$ws1.Range("F18").Value = 626
$ws1.Range("F26").Value = 538
$ws1.Range("F27").Value = 39
$ws1.Range("F31").Value = 3
$ws1.Range("F33").Value = 272
$ws1.Range("F37").Value = 5723
$ws1.Range("F38").Value = 5254

# --- Sheet: 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 80

# --- Sheet: 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 16454
$ws4.Range("F5").Value = 739
$ws4.Range("F6").Value = 15619
$ws4.Range("F8").Value = 9286
$ws4.Range("F9").Value = 500
$ws4.Range("G9").Value = 70
$ws4.Range("G10").Value = 55
$ws4.Range("F11").Value = 1035
$ws4.Range("F12").Value = 130
$ws4.Range("G13").Value = 50
$ws4.Range("G15").Value = 60
$ws4.Range("F16").Value = 23
$ws4.Range("F18").Value = 626
$ws4.Range("F26").Value = 538
$ws4.Range("F27").Value = 39
$ws4.Range("F29").Value = 80
$ws4.Range("F33").Value = 3
$ws4.Range("F35").Value = 272
$ws4.Range("F39").Value = 5723
$ws4.Range("F41").Value = 5254
